# Wenglor relative height uncertainty: decrease from +/- 1 cm (0.01) to
# +/- 1 mm (0.001) in column I (HeightErr_m) for all Wenglor rows
# (rows 506 through 648).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I506:I648").Value = 0.001

# Restore the view: scroll the frozen sheet down near the bottom of the
# data and leave the final selection on M650 (bottom pane).
$ws.Range("M650").Select()
